## Fix Training Data Issue
## The "Date" column (BF) for this sheet held the source filename-derived
## string "5-30-2011-12" instead of the actual game date. Correct every
## data row (BF2:BF31) to the proper ISO date "2012-05-30".
##
## The replacement text looks like a date, so a plain Value2 assignment
## would be auto-converted by Excel into a date serial. To keep it as a
## literal text string (matching the rest of the column), we briefly mark
## the cell as Text (@) before writing it, then restore the cell's
## original (default) formatting by copying the format from the
## unaffected neighboring cell in the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$dateCol = 58   # column BF
$fmtSourceCol = 57  # column BE - same row, untouched default formatting

$oldValue = "5-30-2011-12"
$newValue = "2012-05-30"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)

    if ($cell.Value2 -eq $oldValue) {
        # Force text interpretation so "2012-05-30" isn't coerced to a date serial.
        $cell.NumberFormat = "@"
        $cell.Value2 = $newValue

        # Restore the cell's normal (default) style/format, matching the
        # rest of the sheet, by copying formatting from the adjacent cell.
        $fmtSource = $ws.Cells.Item($r, $fmtSourceCol)
        $fmtSource.Copy()
        $cell.PasteSpecial(-4122)  # xlPasteFormats
    }
}

$excel.CutCopyMode = 0
